$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C11").Value = "yes"
$ws.Range("C12").Value = "yes"
$ws.Range("C13").Value = "yes"
